$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 632.1667
$ws.Range("I18").Value = 700.25
$ws.Range("J18").Value = 496
$ws.Range("K18").Value = 700.25
$ws.Range("L18").Value = 496
$ws.Range("M18").Value = -416.25
$ws.Range("N18").Value = -1064
$ws.Range("H33").Value = 663.1875
$ws.Range("I33").Value = 663.1875
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 663.1875
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -434.1875
$ws.Range("H40").Value = 911.8461
$ws.Range("I40").Value = 745.25
$ws.Range("J40").Value = 985.8889
$ws.Range("K40").Value = 745.25
$ws.Range("L40").Value = 985.8889
$ws.Range("M40").Value = -570.25
$ws.Range("N40").Value = -1335.8889
$ws.Range("H97").Value = 1026.8182
$ws.Range("J97").Value = 1029.5
$ws.Range("L97").Value = 3088.5
$ws.Range("N97").Value = -4080.5
$ws.Range("H112").Value = 4696.904
$ws.Range("J112").Value = 4854.78
$ws.Range("L112").Value = 14564.34
$ws.Range("N112").Value = -16780.34
$ws.Range("H125").Value = 4167686.5
$ws.Range("I125").Value = 12500688
$ws.Range("J125").Value = 1185.75
$ws.Range("K125").Value = 112506192
$ws.Range("L125").Value = 10671.75
$ws.Range("M125").Value = -112503732
$ws.Range("N125").Value = -15591.75
$ws.Range("H127").Value = 777.3889
$ws.Range("I127").Value = 469.7
$ws.Range("J127").Value = 1162
$ws.Range("K127").Value = 1409.1
$ws.Range("L127").Value = 3486
$ws.Range("M127").Value = 3550.9
$ws.Range("N127").Value = -13406
$ws.Range("H129").Value = 1615.9565
$ws.Range("I129").Value = 500
$ws.Range("J129").Value = 1850.8948
$ws.Range("K129").Value = 1500
$ws.Range("L129").Value = 5552.6844
$ws.Range("M129").Value = 3500
$ws.Range("N129").Value = -15552.6844
$ws.Range("H132").Value = 2882.4807
$ws.Range("I132").Value = 2782.5293
$ws.Range("J132").Value = 7980
$ws.Range("K132").Value = 8347.5879
$ws.Range("L132").Value = 23940
$ws.Range("M132").Value = -5817.5879
$ws.Range("N132").Value = -29000
$ws.Range("H133").Value = 65140
$ws.Range("J133").Value = 65140
$ws.Range("L133").Value = 65140
$ws.Range("N133").Value = -75260
$ws.Range("H135").Value = 14707264
$ws.Range("I135").Value = 17242382
$ws.Range("J135").Value = 3579.6
$ws.Range("K135").Value = 155181438
$ws.Range("L135").Value = 32216.4
$ws.Range("M135").Value = -155178903
$ws.Range("N135").Value = -37286.39999999999
$ws.Range("H138").Value = 5284.8
$ws.Range("I138").Value = 8839.4
$ws.Range("J138").Value = 4929.34
$ws.Range("K138").Value = 26518.2
$ws.Range("L138").Value = 14788.02
$ws.Range("M138").Value = -21378.2
$ws.Range("N138").Value = -25068.02
$ws.Range("H140").Value = 76697.14
$ws.Range("J140").Value = 76697.14
$ws.Range("L140").Value = 76697.14
$ws.Range("N140").Value = -87057.14

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 13293.5
$ws.Range("I28").Value = 5848.4
$ws.Range("K28").Value = 5848.4
$ws.Range("M28").Value = -5656.4
$ws.Range("H32").Value = 30325268
$ws.Range("I32").Value = 34503304
$ws.Range("J32").Value = 34499.25
$ws.Range("K32").Value = 34503304
$ws.Range("L32").Value = 34499.25
$ws.Range("M32").Value = -34503017
$ws.Range("N32").Value = -35073.25
$ws.Range("H99").Value = 13293.5
$ws.Range("I99").Value = 5848.4
$ws.Range("K99").Value = 5848.4
$ws.Range("M99").Value = -2853.4
$ws.Range("H141").Value = 44185.6
$ws.Range("J141").Value = 44185.6
$ws.Range("L141").Value = 44185.6
$ws.Range("N141").Value = -54545.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 57646.668
$ws.Range("J51").Value = 57646.668
$ws.Range("L51").Value = 57646.668
$ws.Range("N51").Value = -58628.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 455.91666
$ws.Range("I22").Value = 290.125
$ws.Range("J22").Value = 787.5
$ws.Range("K22").Value = 290.125
$ws.Range("L22").Value = 787.5
$ws.Range("M22").Value = 59.875
$ws.Range("N22").Value = -1487.5
$ws.Range("H31").Value = 11114
$ws.Range("I31").Value = 1262.375
$ws.Range("J31").Value = 13831.689
$ws.Range("K31").Value = 1262.375
$ws.Range("L31").Value = 13831.689
$ws.Range("M31").Value = -967.375
$ws.Range("N31").Value = -14421.689
$ws.Range("H34").Value = 11114
$ws.Range("I34").Value = 1262.375
$ws.Range("J34").Value = 13831.689
$ws.Range("K34").Value = 1262.375
$ws.Range("L34").Value = 13831.689
$ws.Range("M34").Value = -1060.375
$ws.Range("N34").Value = -14235.689
$ws.Range("H50").Value = 19198.8
$ws.Range("J50").Value = 19198.8
$ws.Range("L50").Value = 19198.8
$ws.Range("N50").Value = -20448.8
$ws.Range("H59").Value = 22499
$ws.Range("J59").Value = 22499
$ws.Range("L59").Value = 22499
$ws.Range("N59").Value = -24789
$ws.Range("H68").Value = 25599.5
$ws.Range("J68").Value = 25599.5
$ws.Range("L68").Value = 25599.5
$ws.Range("N68").Value = -27097.5
$ws.Range("H71").Value = 25599.5
$ws.Range("J71").Value = 25599.5
$ws.Range("L71").Value = 76798.5
$ws.Range("N71").Value = -84286.5
$ws.Range("H74").Value = 28888.334
$ws.Range("J74").Value = 28888.334
$ws.Range("L74").Value = 28888.334
$ws.Range("N74").Value = -30636.334
$ws.Range("H77").Value = 28888.334
$ws.Range("J77").Value = 28888.334
$ws.Range("L77").Value = 86665.00199999999
$ws.Range("N77").Value = -95401.00199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1037.473
$ws.Range("I5").Value = 787.1905
$ws.Range("J5").Value = 2470.9092
$ws.Range("K5").Value = 2361.5715
$ws.Range("L5").Value = 7412.7276
$ws.Range("M5").Value = -2249.5715
$ws.Range("N5").Value = -7636.7276
$ws.Range("H81").Value = 5433.2
$ws.Range("J81").Value = 9642.714
$ws.Range("L81").Value = 28928.142
$ws.Range("N81").Value = -31174.142
$ws.Range("H84").Value = 5433.2
$ws.Range("J84").Value = 9642.714
$ws.Range("L84").Value = 86784.42600000001
$ws.Range("N84").Value = -98016.42600000001
$ws.Range("H113").Value = 1037.0454
$ws.Range("I113").Value = 817.16
$ws.Range("K113").Value = 2451.48
$ws.Range("M113").Value = -281.48
$ws.Range("H122").Value = 3057.4038
$ws.Range("I122").Value = 530.7143
$ws.Range("J122").Value = 4769.032
$ws.Range("K122").Value = 4776.428699999999
$ws.Range("L122").Value = 42921.288
$ws.Range("M122").Value = -2326.428699999999
$ws.Range("N122").Value = -47821.288
$ws.Range("H131").Value = 4050.8206
$ws.Range("I131").Value = 351.8
$ws.Range("J131").Value = 4594.794
$ws.Range("K131").Value = 1055.4
$ws.Range("L131").Value = 13784.382
$ws.Range("M131").Value = 3984.6
$ws.Range("N131").Value = -23864.382
$ws.Range("H135").Value = 1037.473
$ws.Range("I135").Value = 787.1905
$ws.Range("J135").Value = 2470.9092
$ws.Range("K135").Value = 7084.7145
$ws.Range("L135").Value = 22238.1828
$ws.Range("M135").Value = -4549.7145
$ws.Range("N135").Value = -27308.1828

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1488.7778
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 1771.2858
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 1771.2858
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -2147.2858
$ws.Range("H133").Value = 49973.57
$ws.Range("J133").Value = 49973.57
$ws.Range("L133").Value = 49973.57
$ws.Range("N133").Value = -55033.57

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 416.12903
$ws.Range("I107").Value = 373.10715
$ws.Range("J107").Value = 817.6667
$ws.Range("K107").Value = 1119.32145
$ws.Range("L107").Value = 2453.0001
$ws.Range("M107").Value = 800.6785500000001
$ws.Range("N107").Value = -6293.0001
$ws.Range("H122").Value = 2774
$ws.Range("I122").Value = 2285.9333
$ws.Range("J122").Value = 3384.0833
$ws.Range("K122").Value = 6857.7999
$ws.Range("L122").Value = 10152.2499
$ws.Range("M122").Value = -4407.7999
$ws.Range("N122").Value = -15052.2499
